$wb = $excel.ActiveWorkbook

# The "想去人数" (number of people interested) column F needs updating
# on both the "展览" sheet and the "全部类型" sheet, which mirror the
# same rows of data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 257
    $ws.Range("F3").Value = 375
}
